# Set explicit column widths on the active sheet so the saved OOXML gets a
# <cols> block (min/max/width/customWidth) matching the target widths:
# A=22, B=20, C=20, D=19, E=20, F=20, G=20, H=20, I=6
#
# Excel stores column width padded by the default character padding
# (5 pixels / Maximum Digit Width, which for the default Calibri 11 font is
# 5/6 of a character), i.e. StoredWidth = Round((ColumnWidth + 5/6) * 6) / 6.
# To land on an exact integer StoredWidth we therefore request
# ColumnWidth = target - 5/6 via the COM ColumnWidth property.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetWidths = @(22, 20, 20, 19, 20, 20, 20, 20, 6)

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $col = $i + 1
    $target = $targetWidths[$i]
    $ws.Columns.Item($col).ColumnWidth = $target - (5.0 / 6.0)
}
